# Update "Training Dashboard" sheet, rows 3-20:
#   - PERIOD TO EXPIRE (col H): decrement numeric value by 1
#   - LAST UPDATE (col I): change text "03-Nov-2025" -> "04-Nov-2025"
#     (written via a formula + paste-values round trip so Excel keeps it
#      as literal text instead of auto-converting it to a date serial)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Scratch cell used to build a "04-Nov-2025" value that Excel will not
# reinterpret as a date when pasted as a value into column I.
$scratch = $ws.Range("Z1")
$scratch.Formula = '="04-Nov-2025"'
$scratch.Copy()

for ($row = 3; $row -le 20; $row++) {
    $hCell = $ws.Cells.Item($row, 8)
    $hCell.Value2 = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)
    $iCell.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

$scratch.ClearContents() | Out-Null
